$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 79. This shifts the existing rows
# 79-108 down to 80-109 (all their values stay intact).
$ws.Rows.Item(79).EntireRow.Insert()

# Populate the newly inserted row 79 with the new weekly price record.
$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 44559
$ws.Range("D79").NumberFormat = $ws.Range("D80").NumberFormat
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100112001
$ws.Cells.Item(79, 7).Value = "Berenjena"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 540
$ws.Cells.Item(79, 11).Value = 9500
$ws.Cells.Item(79, 12).Value = 10000
$ws.Cells.Item(79, 13).Value = 9750
$ws.Cells.Item(79, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 162
$ws.Cells.Item(79, 17).Value = 60
$ws.Cells.Item(79, 18).Value = "Hortaliza"
